# LH_TC_NOTIFICATION.xlsx - v1.1 update
# "Updated the Testcases according to the review"
# LH_TC_NOTIFICATION_Create_028
#
# Adds a new "V1.1" row to the VESRION HISTORY log sheet (with author,
# updated-section note and today's date), and switches the active/visible
# tab from LH_TC_FEATURENAME to VESRION HISTORY.

$wb = $excel.ActiveWorkbook
$wsFeature = $wb.Worksheets.Item("LH_TC_FEATURENAME")
$wsHistory = $wb.Worksheets.Item("VESRION HISTORY")

# --- Append the new version-history row (row 3) ----------------------------
# Copy row 2's look (fonts/fills/borders) down to row 3 first so the new
# row matches the rest of the log, then overwrite with the new content.
$wsHistory.Range("A2:D2").Copy($wsHistory.Range("A3:D3"))

$wsHistory.Range("A3").Value = "V1.1"
$wsHistory.Range("B3").Value = "Mahmoud Abdelmageed"
$wsHistory.Range("C3").Value = "Updated the Testcases according to the review"
$wsHistory.Range("C3").WrapText = $true
$wsHistory.Range("D3").Formula = "=TODAY()"

# Row grows to fit the wrapped "Updated section" text.
$wsHistory.Range("A3:D3").RowHeight = 63

# --- Switch the visible/active sheet to the history log ---------------------
$wsHistory.Activate()
$wsHistory.Range("G14").Select()
